# Add the missing trip-index rows (5 through 13) to column A, continuing
# the existing numeric sequence that ends at row 4 (value 3.0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4.0, 5.0, 7.0, 8.0, 9.0, 10.0, 11.0, 12.0, 13.0)

$row = 5
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}
